$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text-typed numeric-looking value into a cell while
# preserving the cell's original style/number format. Plain assignment of a
# numeric-looking string (e.g. "35.88") auto-converts the cell to a Number,
# but the source data keeps these as text strings, so we briefly force a
# text format, assign, then restore the original style.
function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Enterprises density (per 1000 people) -- row 11
Set-TextValue "B11" "35.88"
Set-TextValue "C11" "8.79"
Set-TextValue "D11" "44.66"

# Employment (% of total) -- row 12
Set-TextValue "B12" "26.98"
Set-TextValue "C12" "60.15"
Set-TextValue "D12" "87.13"

# Enterprises (% of total) -- row 14
Set-TextValue "B14" "80.18"
Set-TextValue "C14" "19.63"
Set-TextValue "D14" "99.81"
